$p = $ppt.ActivePresentation

# --- Slide 3 (Physicians): merge "Patients will " + "only " + "be associated with one physician." into a single run ---
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$para2 = $tr3.Paragraphs(2, 1)
# Force-collapse to a single run by first assigning unrelated text, then the final text.
$para2.Text = "__tmp__"
$para2b = $tr3.Paragraphs(2, 1)
$para2b.Text = "Patients will only be associated with one physician."

# --- Slide 4 (Patients) ---
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: merge "Patient accounts will not " + "have any personal identifiable information" into one run.
$para1 = $tr4.Paragraphs(1, 1)
$para1.Text = "__tmp__"
$para1b = $tr4.Paragraphs(1, 1)
$para1b.Text = "Patient accounts will not have any personal identifiable information"

# Paragraph 2: merge the three runs (PII / call-back email) into a single run.
$para2c = $tr4.Paragraphs(2, 1)
$para2c.Text = "__tmp__"
$para2d = $tr4.Paragraphs(2, 1)
$para2d.Text = "Since there is no PII, password resetting will be performed with either security questions or a temporary password will be mailed to a " + [char]8220 + "call-back email" + [char]8221

# Paragraph 3: surgically insert ", height " before "and location" inside the existing run,
# splitting it into three runs while preserving the untouched "userid" run and its formatting.
$para3 = $tr4.Paragraphs(3, 1)
$para3Text = $para3.Text
$needle = ", and location"
$relIdx = $para3Text.IndexOf($needle)
$absPos = $para3.Start + $relIdx
$commaSpace = $tr4.Characters($absPos, 2)
$commaSpace.Text = ", height "

$ppt.ActivePresentation.Save()
